$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = -12.214
$ws.Range("C12").Value = -12.53
$ws.Range("E13").Value = 12.817
$ws.Range("C18").Value = -12.283
